$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    first (Heading1) paragraph, matching the bold-label + plain
#    text run layout used elsewhere in the document (leading empty
#    run, bold "Meta description" run, plain suffix run).
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$insertPoint = $d.Range($p1.Range.End, $p1.Range.End)

$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Dazzle Me Megaways online slot game with up to 99,225 ways to win and a chance to activate amazing Dazzling Wild Reels for free.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($metaXml)

# InsertXML leaves a stray empty paragraph behind (it was only there
# to force a clean paragraph split) -- remove it.
$d1 = $word.ActiveDocument
$stray = $d1.Paragraphs.Item(3)
$stray.Range.Delete()

# ---------------------------------------------------------------
# 2) Remove the duplicated title paragraph ("Play Dazzle Me
#    Megaways Free | Exciting Wild Reels Feature") that used to sit
#    right before the closing italic description paragraph. Walk
#    paragraphs from the end, skip the genuine Heading1 title (#1),
#    and delete the bold duplicate near the bottom of the document.
# ---------------------------------------------------------------
$d2 = $word.ActiveDocument
$count = $d2.Paragraphs.Count
for ($i = $count; $i -ge 2; $i--) {
    $p = $d2.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Play Dazzle Me Megaways Free | Exciting Wild Reels Feature*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------
# 3) Replace the closing italic description text with the new
#    image-generation prompt text. Scope the Find to the final
#    paragraph only so the newly inserted meta-description
#    paragraph (which shares the same sentence) is left untouched.
# ---------------------------------------------------------------
$d3 = $word.ActiveDocument
$count3 = $d3.Paragraphs.Count
$lastPara = $d3.Paragraphs.Item($count3)
$lastPara.Range.Find.Execute(
    "Read our review of Dazzle Me Megaways online slot game with up to 99,225 ways to win and a chance to activate amazing Dazzling Wild Reels for free.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a feature image for Dazzle Me Megaways that showcases the fun and energetic feel of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses to represent the adventurous and exciting aspect of the slot game. The warrior should be surrounded by colorful gemstones and fruits, which are symbols in the game. The background should have a bright and vibrant effect that brings out the lively experience of playing Dazzle Me Megaways.",
    2)

Write-Host "Done"
